# Commit: "edited condition values from strings to numerical"
#
# The "combined" worksheet has a `condition` column (column E, rows 2-142)
# that was stored as text labels ("DMSO", "WIN05μM", "WIN1μM"). This script
# converts those text labels to their equivalent numeric concentration
# values (0, 0.5, 1 respectively), leaving every other column untouched.
#
# It also updates the worksheet's remembered selection from I10 to K29,
# matching the state the workbook was left in when it was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("combined")

$muChar = [char]0x03BC

for ($r = 2; $r -le 142; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $v = $cell.Value2

    if ($v -eq "DMSO") {
        $cell.Value = 0
    } elseif ($v -eq ("WIN05" + $muChar + "M")) {
        $cell.Value = 0.5
    } elseif ($v -eq ("WIN1" + $muChar + "M")) {
        $cell.Value = 1
    }
}

$ws.Activate()
$ws.Range("K29").Select()
